$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 is already a blank spacer row in the sheet; inserting one row at
# position 42 pushes the old "rc_obs" table (old rows 42-51) down to 43-52,
# turning the old row 42 into a second blank spacer row.
$ws.Rows.Item(42).Insert()

# Populate the new observed-variable row 41 (id / title / title_short / construction).
$ws.Range("A41").Value = "cnds_nom_demean_obs"
$ws.Range("B41").Value = "Nominal cons. on non-durables and services, net growth, demean"
$ws.Range("C41").Value = "Nominal consumption growth"
$ws.Range("D41").Value = "demean:" + [char]0x0394 + "LN(PCEND+PCES)*100"

# Give the "title" cells of this block (B36:B39 existing rows, plus the new
# B41/B42) the same highlighted look already used on B40 (fill, no border).
$ws.Range("B40").Copy()
$top = $ws.Range("B36:B39")
$top.PasteSpecial(-4122)
$top.Borders.LineStyle = -4142

$ws.Range("B40").Copy()
$bottom = $ws.Range("B41:B42")
$bottom.PasteSpecial(-4122)
$bottom.Borders.LineStyle = -4142

$excel.CutCopyMode = 0

# Leave the view scrolled near the edited rows, matching where the author left it.
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("D41").Select()
